$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 22 (Iñigo Arguibide) entirely; row 23 (Sergio Herrera) shifts up to row 22
$ws.Rows("22:22").Delete()
